$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    # Force the cell to Text format so Excel does not reinterpret
    # numeric-looking strings (e.g. "1.01", "4.00") as numbers, then
    # restore General formatting so no visible number-format change
    # is introduced (matches the source workbook, which carries no
    # explicit number formats on these data cells).
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.NumberFormat = "General"
}

Set-TextValue 'D2' '26.129.85'
Set-TextValue 'E2' '  +3.43%  '
Set-TextValue 'D3' '1.598.69'
Set-TextValue 'E3' '  +2.23%  '
Set-TextValue 'D4' '1.01'
Set-TextValue 'E4' '  -0.12%  '
Set-TextValue 'D5' '212.47'
Set-TextValue 'E5' '  +2.65%  '
Set-TextValue 'D6' '1.01'
Set-TextValue 'E6' '  -0.05%  '
Set-TextValue 'D7' '0.484'
Set-TextValue 'E7' '  +2.08%  '
Set-TextValue 'E8' '  +2.60%  '
Set-TextValue 'D9' '0.0616'
Set-TextValue 'E9' '  +1.85%  '
Set-TextValue 'D10' '17.96'
Set-TextValue 'E10' '  +1.25%  '
Set-TextValue 'D11' '0.0821'
Set-TextValue 'E11' '  +4.91%  '
Set-TextValue 'D12' '1.829.67'
Set-TextValue 'E12' '  +2.75%  '
Set-TextValue 'D13' '1.610.24'
Set-TextValue 'E13' '  +2.88%  '
Set-TextValue 'D14' '4.00'
Set-TextValue 'E14' '  -0.47%  '
Set-TextValue 'D15' '0.510'
Set-TextValue 'E15' '  +1.09%  '
Set-TextValue 'D16' '26.106.97'
Set-TextValue 'E16' '  +3.32%  '
Set-TextValue 'D17' '60.59'
Set-TextValue 'E17' '  +2.25%  '
Set-TextValue 'D18' '0.0₃0721'
Set-TextValue 'E18' '  +1.57%  '
Set-TextValue 'B19' 'Dai'
Set-TextValue 'C19' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D19' '1.01'
Set-TextValue 'E19' '  -0.11%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '207.66'
Set-TextValue 'E20' '  +12.06%  '
Set-TextValue 'D21' '4.23'
Set-TextValue 'E21' '  +2.73%  '
Set-TextValue 'D22' '9.31'
Set-TextValue 'E22' '  +0.61%  '
Set-TextValue 'D23' '5.97'
Set-TextValue 'E23' '  +1.80%  '
Set-TextValue 'E24' '  +11.64%  '
Set-TextValue 'D25' '141.74'
Set-TextValue 'E25' '  +1.33%  '
Set-TextValue 'D26' '1.01'
Set-TextValue 'E26' '  -0.17%  '
Set-TextValue 'D27' '0.124'
Set-TextValue 'E27' '  -3.48%  '
Set-TextValue 'E28' '  +2.75%  '
Set-TextValue 'D29' '6.43'
Set-TextValue 'E29' '  +0.00%  '
Set-TextValue 'E30' '  +1.48%  '
Set-TextValue 'D31' '0.0469'
Set-TextValue 'E31' '  +1.55%  '
Set-TextValue 'D32' '3.13'
Set-TextValue 'E32' '  +3.18%  '
Set-TextValue 'D33' '2.99'
Set-TextValue 'E33' '  +0.26%  '
Set-TextValue 'E34' '  +1.31%  '
Set-TextValue 'E35' '  +2.28%  '
Set-TextValue 'D36' '1.109.17'
Set-TextValue 'E36' '  +2.15%  '
Set-TextValue 'D37' '0.0161'
Set-TextValue 'E37' '  +8.14%  '
Set-TextValue 'D38' '1.01'
Set-TextValue 'E38' '  +0.29%  '
Set-TextValue 'D39' '2.33'
Set-TextValue 'E39' '  +0.36%  '
Set-TextValue 'D40' '0.778'
Set-TextValue 'E40' '  +1.28%  '
Set-TextValue 'D41' '0.492'
Set-TextValue 'E41' '  -0.31%  '
Set-TextValue 'D42' '0.777'
Set-TextValue 'E42' '  +1.82%  '
Set-TextValue 'D43' '1.741.62'
Set-TextValue 'E43' '  +2.70%  '
Set-TextValue 'D44' '92.63'
Set-TextValue 'E44' '  -0.31%  '
Set-TextValue 'E46' '  -0.31%  '
Set-TextValue 'E47' '  +6.76%  '
Set-TextValue 'D48' '53.49'
Set-TextValue 'E48' '  +1.61%  '
Set-TextValue 'E49' '  +0.24%  '
Set-TextValue 'E50' '  +0.82%  '
Set-TextValue 'E51' '  +0.08%  '
